$p = $ppt.ActivePresentation

# 1. Handout master date placeholder: 2024-10-20 -> 2024-10-31
$hm = $p.HandoutMaster
$hmDate = $hm.Shapes.Item(2)
$hmDate.TextFrame.TextRange.Text = "2024-10-31"

# 2. Notes master date placeholder (fr-FR): 20/10/2024 -> 31/10/2024
$nm = $p.NotesMaster
$nmDate = $nm.Shapes.Item(2)
$nmDate.TextFrame.TextRange.Text = "31/10/2024"

# 3. Slide 1 "TextBox 1": add workshop lead name and date
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$tr.Paragraphs(1, 1).Text = "Benjamin Rudski"
$tr.Paragraphs(2, 1).Text = "November 1, 2024"
